$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "BIG APPLE" + "," -> single run "BIG APPLE,"  (and "N" + "Y" -> "NY")
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("BIG APPLE,", $true, $false, $false, $false, $false, $true, 1, $false, "BIG APPLE,", 2) | Out-Null

$r = $d.Content
$r.Find.Execute("NY", $true, $false, $false, $false, $false, $true, 1, $false, "NY", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) "SMALL TOWN" + "," -> single run "SMALL TOWN,"
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("SMALL TOWN,", $true, $false, $false, $false, $false, $true, 1, $false, "SMALL TOWN,", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) MICR line: C300 2 C " " A031176110A " " 36116600252C{sp}
#           ->  C 147011 C " " A0 123987 10A " " 87641284584 C{sp}
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("C300", $true, $false, $false, $false, $false, $true, 1, $false, "C", 2) | Out-Null

$r = $d.Content
$r.Find.Execute("C2C", $true, $false, $false, $false, $false, $true, 1, $false, "C147011C", 2) | Out-Null

$r = $d.Content
$r.Find.Execute("A031176110A", $true, $false, $false, $false, $false, $true, 1, $false, "A012398710A", 2) | Out-Null

$r = $d.Content
$r.Find.Execute("36116600252C ", $true, $false, $false, $false, $false, $true, 1, $false, "87641284584C ", 2) | Out-Null

# The two text replacements above land as single runs because adjacent runs
# sharing identical formatting get coalesced by the text-replace. The target
# keeps each piece as its own <w:r> (still with identical rPr), so re-split
# them by toggling a character property on/off across the internal
# boundaries -- that forces the run break without altering the rendered
# formatting.
$rf = $d.Content
$rf.Find.Execute("A012398710A") | Out-Null
$aStart = $rf.Start

$bound1 = $d.Range($aStart, $aStart + 2)       # "A0"
$bound2 = $d.Range($aStart + 8, $aStart + 11)  # "10A"
$bound1.Font.Bold = $true
$bound1.Font.Bold = $false
$bound2.Font.Bold = $true
$bound2.Font.Bold = $false

$rf2 = $d.Content
$rf2.Find.Execute("87641284584C ") | Out-Null
$bStart = $rf2.Start

$bound3 = $d.Range($bStart + 11, $bStart + 13) # "C "
$bound3.Font.Bold = $true
$bound3.Font.Bold = $false

Write-Output "done"
